{"js": "// 1) UC3.2 subcase: \"Visualizzazione buoni in circolazione\" ->\n//    \"Eliminazione automatica buoni sconto compleanno scaduti\"\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet ucPara = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"Visualizzazione buoni in circolazione\") !== -1) {\n    ucPara = p;\n    break;\n  }\n}\n\nif (ucPara) {\n  // Scope the search to this paragraph only, so we don't touch the similarly\n  // worded \"UC1.2: Visualizzazione appuntamenti\" elsewhere in the document.\n  const hit = ucPara.search(\"Visualizzazione buoni in circolazione\", { matchCase: true });\n  hit.load(\"text\");\n  await context.sync();\n  if (hit.items.length > 0) {\n    hit.items[0].insertText(\n      \"Eliminazione automatica buoni sconto compleanno scaduti\",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n}\n\n// 2) Reorder the summary table: move the \"UC3 / Gestione buoni sconto\n//    compleanno\" row so it comes right after \"UC1 / Gestione appuntamenti\"\n//    (i.e. before \"UC2 / Gestione clienti\") instead of being last.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length > 0) {\n  const table = tables.items[0];\n  table.rows.load(\"items\");\n  await context.sync();\n\n  const rows = table.rows.items;\n  // Read every row's two cell texts.\n  for (const row of rows) {\n    row.cells.load(\"items\");\n  }\n  await context.sync();\n  for (const row of rows) {\n    for (const cell of row.cells.items) {\n      cell.body.load(\"text\");\n    }\n  }\n  await context.sync();\n\n  let uc1Row = null;\n  let uc3Row = null;\n  let uc3Code = \"UC3\";\n  let uc3Title = \"Gestione buoni sconto compleanno\";\n  for (const row of rows) {\n    const texts = row.cells.items.map((c) => c.body.text.trim());\n    if (texts[0] === \"UC1\") {\n      uc1Row = row;\n    } else if (texts[0] === \"UC3\") {\n      uc3Row = row;\n      uc3Code = texts[0];\n      uc3Title = texts[1];\n    }\n  }\n\n  if (uc1Row && uc3Row) {\n    // Insert a fresh row (single run per cell, just like the target OOXML)\n    // right after the UC1 row, carrying the UC3 content.\n    uc1Row.insertRows(Word.InsertLocation.after, 1, [[uc3Code, uc3Title]]);\n    await context.sync();\n\n    // Re-fetch the table rows: after the insertion the previously captured\n    // row objects/indices are stale, so look the old (now duplicate),\n    // trailing UC3 row up again before removing it.\n    const tables2 = context.document.body.tables;\n    tables2.load(\"items\");\n    await context.sync();\n    const table2 = tables2.items[0];\n    table2.rows.load(\"items\");\n    await context.sync();\n    const rows2 = table2.rows.items;\n    for (const row of rows2) {\n      row.cells.load(\"items\");\n    }\n    await context.sync();\n    for (const row of rows2) {\n      for (const cell of row.cells.items) {\n        cell.body.load(\"text\");\n      }\n    }\n    await context.sync();\n\n    let uc3Rows = [];\n    for (const row of rows2) {\n      const texts = row.cells.items.map((c) => c.body.text.trim());\n      if (texts[0] === \"UC3\") {\n        uc3Rows.push(row);\n      }\n    }\n    // The duplicate (original) UC3 row is now the last one in the table.\n    if (uc3Rows.length > 1) {\n      uc3Rows[uc3Rows.length - 1].delete();\n      await context.sync();\n    }\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) UC3.2 subcase: \"Visualizzazione buoni in circolazione\" ->\n#    \"Eliminazione automatica buoni sconto compleanno scaduti\"\n# Locate the exact paragraph first so we don't touch the similarly worded\n# \"UC1.2: Visualizzazione appuntamenti\" earlier in the document.\nforeach ($p in $d.Paragraphs) {\n    $ptext = $p.Range.Text\n    if ($ptext -like \"*Visualizzazione buoni in circolazione*\") {\n        $rng = $p.Range\n        $rng.Find.Execute(\"Visualizzazione buoni in circolazione\", $false, $false, $false, $false, $false, $true, 1, $false, \"Eliminazione automatica buoni sconto compleanno scaduti\", 2)\n        break\n    }\n}\n\n# 2) Reorder the summary table: move the \"UC3 / Gestione buoni sconto\n#    compleanno\" row so it comes right after \"UC1 / Gestione appuntamenti\"\n#    (i.e. before \"UC2 / Gestione clienti\") instead of being last.\n$t = $d.Tables.Item(1)\n\n$uc1Index = 0\n$uc2Index = 0\n$uc3Index = 0\n$uc3Code = \"UC3\"\n$uc3Title = \"Gestione buoni sconto compleanno\"\n\nfor ($i = 1; $i -le $t.Rows.Count; $i++) {\n    $row = $t.Rows.Item($i)\n    $code = $row.Cells.Item(1).Range.Text.TrimEnd([char]7, [char]13)\n    if ($code -eq \"UC1\") {\n        $uc1Index = $i\n    } elseif ($code -eq \"UC2\") {\n        $uc2Index = $i\n    } elseif ($code -eq \"UC3\") {\n        $uc3Index = $i\n        $uc3Code = $code\n        $uc3Title = $row.Cells.Item(2).Range.Text.TrimEnd([char]7, [char]13)\n    }\n}\n\nif (($uc1Index -gt 0) -and ($uc3Index -gt 0)) {\n    # Insert a fresh row right before the UC2 row (i.e. right after UC1),\n    # carrying the UC3 content, then drop the old trailing UC3 row.\n    $beforeRow = $t.Rows.Item($uc1Index + 1)\n    $newRow = $t.Rows.Add($beforeRow)\n    $newRow.Cells.Item(1).Range.Text = $uc3Code\n    $newRow.Cells.Item(2).Range.Text = $uc3Title\n\n    # The old UC3 row is now the last row of the table.\n    $lastRow = $t.Rows.Item($t.Rows.Count)\n    $lastRow.Delete()\n}\n"}
